# Add the four "near failure" header columns (K1:N1) to the gate upload
# template, matching the author's "added near failure function" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a transient named style so the new DengXian font is registered
# once (cell-level Font.Name assignment clones/duplicates font entries).
$styleName = "NearFailureHeader"
$style = $wb.Styles.Add($styleName)
$style.Font.Name = "DengXian"
$style.Font.Size = 12

$headers = @("hand_near_max", "hand_near_min", "foot_near_max", "foot_near_min")
$cols = @("K", "L", "M", "N")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Style = $styleName
    $cell.VerticalAlignment = -4108  # xlCenter
}

$wb.Styles($styleName).Delete()
